$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the previous used range entirely so stale columns (K:V) and any
# leftover formatting/values from the old layout are gone.
$ws.Range("A1:V2").Clear()

# New header row
$ws.Range("A1").Value = "_id"
$ws.Range("B1").Value = "fullname"
$ws.Range("C1").Value = "username"
$ws.Range("D1").Value = "email"
$ws.Range("E1").Value = "accountrole"
$ws.Range("F1").Value = "__v"
$ws.Range("G1").Value = "address"
$ws.Range("H1").Value = "birthday"
$ws.Range("I1").Value = "phone"
$ws.Range("J1").Value = "password"

# Row 2
$ws.Range("A2").Value = "62557f02b9324784a2aa21f8"
$ws.Range("B2").Value = "Park Roseann             "
$ws.Range("C2").Value = "roseann"
$ws.Range("D2").Value = "roseann@g.c"
$ws.Range("E2").Value = "admin"
$ws.Range("F2").Value = 0
$ws.Range("G2").Value = "South Korea"
$ws.Range("H2").NumberFormat = "@"
$ws.Range("H2").Value = "02/16/1997"
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "59454"
$ws.Range("J2").NumberFormat = "@"
$ws.Range("J2").Value = "123"

# Row 3
$ws.Range("A3").Value = "6255810366248bbde36e8563"
$ws.Range("B3").Value = "kun"
$ws.Range("C3").Value = "kunx"
$ws.Range("D3").Value = "kunx@g.c"
$ws.Range("E3").Value = "employee"
$ws.Range("F3").Value = 0

# Row 4
$ws.Range("A4").Value = "625582820fec137040de7e78"
$ws.Range("B4").Value = "   adminx   "
$ws.Range("C4").Value = "adminx"
$ws.Range("D4").Value = "adminx@g.c"
$ws.Range("E4").Value = "admin"
$ws.Range("F4").Value = 0
$ws.Range("J4").NumberFormat = "@"
$ws.Range("J4").Value = "234"
